$d = $word.ActiveDocument

# --- 1. Split the first paragraph's run into "Hey there it's me " + proofErr-wrapped "vj" ---
$p1 = $d.Paragraphs(1).Range
$p1.MoveEnd(1, -1)  # exclude the paragraph mark, keep only the run text

$splitXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Hey there it’s me </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>vj</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@
$p1.InsertXML($splitXml)

# --- 2. Append the four new paragraphs at the end of the document ---
# Each InsertXML call is issued against a *freshly fetched* end-of-story
# range, since a previously-used Range object does not auto-advance past
# content it just inserted.

$p2Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I made a change</w:t></w:r></w:p>
"@
$rng2 = $d.Content
$rng2.Collapse(0)
$rng2.InsertXML($p2Xml)

$p3Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>I made a change</w:t></w:r></w:p>
"@
$rng3 = $d.Content
$rng3.Collapse(0)
$rng3.InsertXML($p3Xml)

$p4Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
"@
$rng4 = $d.Content
$rng4.Collapse(0)
$rng4.InsertXML($p4Xml)

$p5Xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>
"@
$rng5 = $d.Content
$rng5.Collapse(0)
$rng5.InsertXML($p5Xml)
